# Fruta / hortaliza, semanal
# Insert a new weekly record as row 14, shifting existing rows 14-35 down to 15-36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14 (entire row), pushing rows 14..35 down to 15..36.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R keep the same constant values used throughout
# this sheet; D,J,K,L,M,P carry the new observation's data.
$ws.Range("A14").Value = 9
$ws.Range("B14").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C14").Value = "Metropolitana"
$ws.Range("D14").Value = 44952
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 100112010
$ws.Range("G14").Value = "Achicoria"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 70
$ws.Range("K14").Value = 7000
$ws.Range("L14").Value = 7000
$ws.Range("M14").Value = 7000
$ws.Range("N14").Value = "`$/caja 16 unidades"
$ws.Range("O14").Value = "Provincia de Quillota"
$ws.Range("P14").Value = 438
$ws.Range("Q14").Value = 16
$ws.Range("R14").Value = "Hortaliza"

# Apply the same date number format used by the other date cells in column D.
$ws.Range("D14").NumberFormat = $ws.Range("D15").NumberFormat
